# Daily attendance processing - reorders "Recorded By" (column G) values so
# that an exact, case-sensitive "System" token is moved to the front of the
# comma-separated list, preserving the relative order of the remaining
# tokens (including a lowercase "system" if present).

function Test-ExactSystem($s) {
    if ($s.Length -ne 6) { return $false }
    $expected = @(83, 121, 115, 116, 101, 109)  # ASCII codes for "System"
    for ($i = 0; $i -lt 6; $i++) {
        $code = [int][char]$s.Substring($i, 1)
        if ($code -ne $expected[$i]) { return $false }
    }
    return $true
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $v = $cell.Value2

    if ($v -eq $null) { continue }
    if ($v.Length -eq 0) { continue }
    if ($v.IndexOf(",") -lt 0) { continue }

    $parts = $v.Split(",")

    $rest = New-Object System.Collections.ArrayList
    $foundIndex = -1
    $idx = 0
    foreach ($p in $parts) {
        $t = $p.Trim()
        if ((Test-ExactSystem $t) -and $foundIndex -eq -1) {
            $foundIndex = $idx
        } else {
            [void]$rest.Add($t)
        }
        $idx = $idx + 1
    }

    if ($foundIndex -gt 0) {
        [void]$rest.Insert(0, "System")
        $result = [string]::Join(", ", $rest)
        $cell.Value2 = $result
    }
}
